# Update the cryptos list (Price and Volume(1h) columns) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "28.284.03"; E = "  +2.66%  " },
    @{ Row = 3; D = "1.871.31"; E = "  +1.44%  " },
    @{ Row = 4; D = $null; E = "  -0.34%  " },
    @{ Row = 5; D = "337.56"; E = "  +1.26%  " },
    @{ Row = 6; D = $null; E = "  -0.32%  " },
    @{ Row = 7; D = $null; E = "  +1.39%  " },
    @{ Row = 8; D = "0.3934"; E = "  +2.04%  " },
    @{ Row = 9; D = "47.23"; E = "  +2.19%  " },
    @{ Row = 10; D = "0.08001"; E = "  +1.14%  " },
    @{ Row = 11; D = "1.005"; E = "  +1.11%  " },
    @{ Row = 12; D = "21.77"; E = "  +1.50%  " },
    @{ Row = 13; D = "1.884.87"; E = "  +2.00%  " },
    @{ Row = 14; D = "5.993"; E = "  +1.29%  " },
    @{ Row = 15; D = "7.290"; E = "  +2.55%  " },
    @{ Row = 16; D = "91.31"; E = "  +2.73%  " },
    @{ Row = 17; D = $null; E = "  -0.39%  " },
    @{ Row = 18; D = "0.00001043"; E = "  +0.86%  " },
    @{ Row = 19; D = "0.06582"; E = "  -0.97%  " },
    @{ Row = 20; D = "17.68"; E = "  +3.70%  " },
    @{ Row = 21; D = $null; E = "  -0.26%  " },
    @{ Row = 22; D = "28.282.93"; E = "  +2.59%  " },
    @{ Row = 23; D = "5.451"; E = "  +1.30%  " },
    @{ Row = 24; D = $null; E = "  +1.27%  " },
    @{ Row = 25; D = "2.295"; E = "  -0.29%  " },
    @{ Row = 26; D = "2.097.70"; E = "  +1.51%  " },
    @{ Row = 27; D = "159.32"; E = "  +0.78%  " },
    @{ Row = 28; D = "19.83"; E = "  +1.78%  " },
    @{ Row = 29; D = "2.159"; E = "  +2.85%  " },
    @{ Row = 30; D = "5.501"; E = "  +1.89%  " },
    @{ Row = 31; D = "119.98"; E = "  +0.12%  " },
    @{ Row = 32; D = "0.9806"; E = "  +0.43%  " },
    @{ Row = 33; D = "0.09501"; E = "  +0.99%  " },
    @{ Row = 34; D = "3.581"; E = "  -0.14%  " },
    @{ Row = 35; D = "1.381"; E = "  +3.20%  " },
    @{ Row = 36; D = "5.361"; E = "  +1.40%  " },
    @{ Row = 37; D = "0.02276"; E = "  +2.18%  " },
    @{ Row = 38; D = "0.06090"; E = "  +1.03%  " },
    @{ Row = 39; D = "8.454"; E = "  +1.66%  " },
    @{ Row = 40; D = "1.178"; E = "  -0.36%  " },
    @{ Row = 41; D = "0.5968"; E = "  +1.39%  " },
    @{ Row = 42; D = $null; E = "  -0.28%  " },
    @{ Row = 43; D = "0.1877"; E = "  +0.76%  " },
    @{ Row = 44; D = "10.40"; E = "  +0.90%  " },
    @{ Row = 45; D = $null; E = "  +4.83%  " },
    @{ Row = 46; D = "0.5620"; E = "  +0.75%  " },
    @{ Row = 47; D = "12.15"; E = "  +0.21%  " },
    @{ Row = 48; D = $null; E = "  +3.51%  " },
    @{ Row = 49; D = "0.06898"; E = "  +3.17%  " },
    @{ Row = 50; D = "110.73"; E = "  -0.08%  " },
    @{ Row = 51; D = "2.020"; E = "  +13.55%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        # Force numeric-looking strings to stay text, matching the source data
        # (e.g. "337.56", "0.3934") without leaving stray formatting behind.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
